# The workbook had one header/data row per line starting at row 1 (A1:E39).
# A new blank row was inserted above the existing data, pushing everything
# down by one row (new used range A2:E40) and leaving row 1 empty.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").EntireRow.Insert()

# Match the saved cursor position recorded in the workbook after the edit.
$ws.Range("I10").Select()
